$d = $word.ActiveDocument
$d.Content.Find.Execute("91×36=3276", $true, $false, $false, $false, $false, $true, 1, $false, "26×28=728", 2)
$d.Content.Find.Execute("57×60=3420", $true, $false, $false, $false, $false, $true, 1, $false, "40×60=2400", 2)
$d.Content.Find.Execute("41×23=943", $true, $false, $false, $false, $false, $true, 1, $false, "88×29=2552", 2)
$d.Content.Find.Execute("79×84=6636", $true, $false, $false, $false, $false, $true, 1, $false, "18×67=1206", 2)
$d.Content.Find.Execute("91×91=8281", $true, $false, $false, $false, $false, $true, 1, $false, "41×67=2747", 2)
$d.Content.Find.Execute("92×29=2668", $true, $false, $false, $false, $false, $true, 1, $false, "90×39=3510", 2)
$d.Content.Find.Execute("11×42=462", $true, $false, $false, $false, $false, $true, 1, $false, "13×61=793", 2)
$d.Content.Find.Execute("82×47=3854", $true, $false, $false, $false, $false, $true, 1, $false, "18×32=576", 2)
$d.Content.Find.Execute("57×84=4788", $true, $false, $false, $false, $false, $true, 1, $false, "58×64=3712", 2)
$d.Content.Find.Execute("46×13=598", $true, $false, $false, $false, $false, $true, 1, $false, "94×42=3948", 2)
$d.Content.Find.Execute("91×88=8008", $true, $false, $false, $false, $false, $true, 1, $false, "69×75=5175", 2)
$d.Content.Find.Execute("54×77=4158", $true, $false, $false, $false, $false, $true, 1, $false, "28×31=868", 2)
$d.Content.Find.Execute("88×68=5984", $true, $false, $false, $false, $false, $true, 1, $false, "22×82=1804", 2)
$d.Content.Find.Execute("82×87=7134", $true, $false, $false, $false, $false, $true, 1, $false, "55×38=2090", 2)
$d.Content.Find.Execute("15×54=810", $true, $false, $false, $false, $false, $true, 1, $false, "79×33=2607", 2)
$d.Content.Find.Execute("43×55=2365", $true, $false, $false, $false, $false, $true, 1, $false, "52×16=832", 2)
$d.Content.Find.Execute("30×93=2790", $true, $false, $false, $false, $false, $true, 1, $false, "16×99=1584", 2)
$d.Content.Find.Execute("56×49=2744", $true, $false, $false, $false, $false, $true, 1, $false, "62×68=4216", 2)
$d.Content.Find.Execute("34×11=374", $true, $false, $false, $false, $false, $true, 1, $false, "46×16=736", 2)
$d.Content.Find.Execute("54×32=1728", $true, $false, $false, $false, $false, $true, 1, $false, "20×95=1900", 2)
$d.Content.Find.Execute("29×78=2262", $true, $false, $false, $false, $false, $true, 1, $false, "52×62=3224", 2)
$d.Content.Find.Execute("71×61=4331", $true, $false, $false, $false, $false, $true, 1, $false, "88×94=8272", 2)
$d.Content.Find.Execute("92×67=6164", $true, $false, $false, $false, $false, $true, 1, $false, "17×82=1394", 2)
$d.Content.Find.Execute("86×22=1892", $true, $false, $false, $false, $false, $true, 1, $false, "59×71=4189", 2)
$d.Content.Find.Execute("77×82=6314", $true, $false, $false, $false, $false, $true, 1, $false, "25×66=1650", 2)
Write-Output "done"
